$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 43.25
$ws.Range("I9").Value = 35.142857
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 35.142857
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 133.857143
$ws.Range("N9").Value = -438
$ws.Range("H12").Value = 312.14285
$ws.Range("I12").Value = 314.23077
$ws.Range("K12").Value = 314.23077
$ws.Range("M12").Value = -144.23077
$ws.Range("H31").Value = 578.8
$ws.Range("I31").Value = 578.8
$ws.Range("K31").Value = 1736.4
$ws.Range("M31").Value = -1506.4
$ws.Range("H33").Value = 39735.79
$ws.Range("J33").Value = 493.4
$ws.Range("L33").Value = 493.4
$ws.Range("N33").Value = -951.4
$ws.Range("H62").Value = 27521.773
$ws.Range("J62").Value = 41317.07
$ws.Range("L62").Value = 41317.07
$ws.Range("N62").Value = -42565.07
$ws.Range("H65").Value = 27521.773
$ws.Range("J65").Value = 41317.07
$ws.Range("L65").Value = 206585.35
$ws.Range("N65").Value = -212825.35
$ws.Range("H92").Value = 73968.11
$ws.Range("I92").Value = 271.73914
$ws.Range("K92").Value = 271.73914
$ws.Range("M92").Value = 976.26086
$ws.Range("H129").Value = 888.9524
$ws.Range("J129").Value = 1311
$ws.Range("L129").Value = 3933
$ws.Range("N129").Value = -13933
$ws.Range("H132").Value = 16669069
$ws.Range("I132").Value = 18520608
$ws.Range("K132").Value = 55561824
$ws.Range("M132").Value = -55559294
$ws.Range("H137").Value = 2239.9644
$ws.Range("I137").Value = 1806.9524
$ws.Range("K137").Value = 5420.857199999999
$ws.Range("M137").Value = -2870.857199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 644.9
$ws.Range("I2").Value = 647.2632
$ws.Range("K2").Value = 647.2632
$ws.Range("M2").Value = -534.2632
$ws.Range("H61").Value = 3366.4167
$ws.Range("I61").Value = 3039.7
$ws.Range("K61").Value = 3039.7
$ws.Range("M61").Value = -2827.7
$ws.Range("H63").Value = 120006350
$ws.Range("I63").Value = 500001760
$ws.Range("J63").Value = 25007500
$ws.Range("K63").Value = 500001760
$ws.Range("L63").Value = 25007500
$ws.Range("M63").Value = -500001074
$ws.Range("N63").Value = -25008872
$ws.Range("H66").Value = 120006350
$ws.Range("I66").Value = 500001760
$ws.Range("J66").Value = 25007500
$ws.Range("K66").Value = 2500008800
$ws.Range("L66").Value = 125037500
$ws.Range("M66").Value = -2500005368
$ws.Range("N66").Value = -125044364
$ws.Range("H74").Value = 40002372
$ws.Range("I74").Value = 47620616
$ws.Range("J74").Value = 6598.5
$ws.Range("K74").Value = 47620616
$ws.Range("L74").Value = 6598.5
$ws.Range("M74").Value = -47619742
$ws.Range("N74").Value = -8346.5
$ws.Range("H77").Value = 40002372
$ws.Range("I77").Value = 47620616
$ws.Range("J77").Value = 6598.5
$ws.Range("K77").Value = 238103080
$ws.Range("L77").Value = 32992.5
$ws.Range("M77").Value = -238098712
$ws.Range("N77").Value = -41728.5
$ws.Range("H110").Value = 28572918
$ws.Range("I110").Value = 37038320
$ws.Range("K110").Value = 37038320
$ws.Range("M110").Value = -37036275
$ws.Range("H116").Value = 644.9
$ws.Range("I116").Value = 647.2632
$ws.Range("K116").Value = 647.2632
$ws.Range("M116").Value = 1646.7368
$ws.Range("H132").Value = 2868.9092
$ws.Range("I132").Value = 1853.8823
$ws.Range("J132").Value = 6320
$ws.Range("K132").Value = 5561.6469
$ws.Range("L132").Value = 18960
$ws.Range("M132").Value = -3031.6469
$ws.Range("N132").Value = -24020
$ws.Range("H136").Value = 3366.4167
$ws.Range("I136").Value = 3039.7
$ws.Range("K136").Value = 9119.099999999999
$ws.Range("M136").Value = -6569.099999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 644.9
$ws.Range("I3").Value = 647.2632
$ws.Range("K3").Value = 647.2632
$ws.Range("M3").Value = -533.2632
$ws.Range("H14").Value = 300
$ws.Range("I14").Value = 300
$ws.Range("K14").Value = 300
$ws.Range("M14").Value = -128
$ws.Range("H36").Value = 5200
$ws.Range("I36").Value = 5200
$ws.Range("K36").Value = 5200
$ws.Range("M36").Value = -4666
$ws.Range("H94").Value = 8622323
$ws.Range("I94").Value = 12501389
$ws.Range("K94").Value = 12501389
$ws.Range("M94").Value = -12500938
$ws.Range("H107").Value = 35720012
$ws.Range("I107").Value = 2327.8572
$ws.Range("J107").Value = 71437700
$ws.Range("K107").Value = 2327.8572
$ws.Range("L107").Value = 71437700
$ws.Range("M107").Value = -407.8571999999999
$ws.Range("N107").Value = -71441540

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1611.6666
$ws.Range("I16").Value = 1567.25
$ws.Range("K16").Value = 1567.25
$ws.Range("M16").Value = -1280.25
$ws.Range("H31").Value = 2160.7104
$ws.Range("I31").Value = 1643.3182
$ws.Range("K31").Value = 1643.3182
$ws.Range("M31").Value = -1348.3182
$ws.Range("H34").Value = 2160.7104
$ws.Range("I34").Value = 1643.3182
$ws.Range("K34").Value = 1643.3182
$ws.Range("M34").Value = -1441.3182
$ws.Range("H113").Value = 1611.6666
$ws.Range("I113").Value = 1567.25
$ws.Range("K113").Value = 1567.25
$ws.Range("M113").Value = 602.75
$ws.Range("H132").Value = 386916.5
$ws.Range("I132").Value = 1978.3
$ws.Range("K132").Value = 5934.9
$ws.Range("M132").Value = -3404.9
$ws.Range("H134").Value = 3467.647
$ws.Range("I134").Value = 2857.1628
$ws.Range("K134").Value = 8571.4884
$ws.Range("M134").Value = -6036.4884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 269
$ws.Range("I5").Value = 269
$ws.Range("K5").Value = 807
$ws.Range("M5").Value = -695
$ws.Range("H135").Value = 269
$ws.Range("I135").Value = 269
$ws.Range("K135").Value = 2421
$ws.Range("M135").Value = 114

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2503237.5
$ws.Range("I7").Value = 4316.6665
$ws.Range("J7").Value = 10000000
$ws.Range("K7").Value = 4316.6665
$ws.Range("L7").Value = 10000000
$ws.Range("M7").Value = -4204.6665
$ws.Range("N7").Value = -10000224
$ws.Range("H8").Value = 2503237.5
$ws.Range("I8").Value = 4316.6665
$ws.Range("J8").Value = 10000000
$ws.Range("K8").Value = 4316.6665
$ws.Range("L8").Value = 10000000
$ws.Range("M8").Value = -4177.6665
$ws.Range("N8").Value = -10000278
$ws.Range("H80").Value = 6463.9375
$ws.Range("J80").Value = 6917.769
$ws.Range("L80").Value = 6917.769
$ws.Range("N80").Value = -8913.769
$ws.Range("H83").Value = 6463.9375
$ws.Range("J83").Value = 6917.769
$ws.Range("L83").Value = 34588.845
$ws.Range("N83").Value = -44572.845
$ws.Range("H113").Value = 3992.611
$ws.Range("J113").Value = 4625.4546
$ws.Range("L113").Value = 4625.4546
$ws.Range("N113").Value = -8965.454600000001
$ws.Range("H132").Value = 3188.2163
$ws.Range("I132").Value = 2837.3044
$ws.Range("J132").Value = 3764.7144
$ws.Range("K132").Value = 8511.913199999999
$ws.Range("L132").Value = 11294.1432
$ws.Range("M132").Value = -5981.913199999999
$ws.Range("N132").Value = -16354.1432
$ws.Range("H136").Value = 27565.77
$ws.Range("J136").Value = 27565.77
$ws.Range("L136").Value = 82697.31
$ws.Range("N136").Value = -87797.31

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3870.65
$ws.Range("I122").Value = 3780.8667
$ws.Range("K122").Value = 11342.6001
$ws.Range("M122").Value = -8892.6001
$ws.Range("H136").Value = 5696.9062
$ws.Range("I136").Value = 5696.72
$ws.Range("J136").Value = 5697.5713
$ws.Range("K136").Value = 17090.16
$ws.Range("L136").Value = 17092.7139
$ws.Range("M136").Value = -14540.16
$ws.Range("N136").Value = -22192.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 505
$ws.Range("I13").Value = 505
$ws.Range("K13").Value = 505
$ws.Range("M13").Value = -365
$ws.Range("H74").Value = 49219.25
$ws.Range("I74").Value = 49625
$ws.Range("J74").Value = 49084
$ws.Range("K74").Value = 49625
$ws.Range("L74").Value = 49084
$ws.Range("M74").Value = -48689
$ws.Range("N74").Value = -50956
$ws.Range("H77").Value = 49219.25
$ws.Range("I77").Value = 49625
$ws.Range("J77").Value = 49084
$ws.Range("K77").Value = 148875
$ws.Range("L77").Value = 147252
$ws.Range("M77").Value = -144195
$ws.Range("N77").Value = -156612
$ws.Range("H80").Value = 35000
$ws.Range("J80").Value = 35000
$ws.Range("L80").Value = 35000
$ws.Range("N80").Value = -36996
$ws.Range("H83").Value = 35000
$ws.Range("J83").Value = 35000
$ws.Range("L83").Value = 105000
$ws.Range("N83").Value = -114984
$ws.Range("H107").Value = 519.3043
$ws.Range("I107").Value = 652.0833
$ws.Range("K107").Value = 1956.2499
$ws.Range("M107").Value = -36.24990000000003
$ws.Range("H113").Value = 1606.0264
$ws.Range("J113").Value = 2364.4211
$ws.Range("L113").Value = 7093.263300000001
$ws.Range("N113").Value = -11433.2633
$ws.Range("H132").Value = 327314.66
$ws.Range("I132").Value = 448700.25
$ws.Range("K132").Value = 1346100.75
$ws.Range("M132").Value = -1343570.75
